$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.219.07'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.630.77'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.89'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.521'
$ws.Range('E6').Value = '  +1.97%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.255'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.21'
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '1.635.28'
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.543'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '27.208.02'
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.88'
$ws.Range('E16').Value = '  -4.22%  '
$ws.Range('D17').Value = '0.0₃0733'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '216.30'
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.94'
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.46'
$ws.Range('E22').Value = '  -3.39%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.11'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '148.07'
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.27'
$ws.Range('E26').Value = '  -3.79%  '
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.57'
$ws.Range('E28').Value = '  -1.64%  '
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.38'
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.00'
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('D33').Value = '1.317.96'
$ws.Range('E33').Value = '  +4.69%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.55'
$ws.Range('E34').Value = '  -2.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.46'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  -1.83%  '
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('E41').Value = '  -1.27%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '63.51'
$ws.Range('E42').Value = '  +2.28%  '
$ws.Range('D43').Value = '1.768.36'
$ws.Range('E43').Value = '  -1.51%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.20'
$ws.Range('E44').Value = '  -4.68%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '90.80'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.820'
$ws.Range('E48').Value = '  +22.19%  '
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.56'
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0957'
$ws.Range('E51').Value = '  -1.67%  '
